# Edit script: reshape "Organizations" sheet to new column layout
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the obsolete "Tiktok Link" column (M) entirely.
$ws.Columns.Item(13).Delete()

# 2. Remove the now-duplicate rows (old rows 3 and 4 duplicated row 2's
#    "Student Life" entry) plus the two extra trailing rows that get
#    folded into the remaining 4 data rows below.
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()

# After the deletes above, sheet has 5 rows left:
#   Row1 = headers, Row2 = Student Life, Row3 = Clubs,
#   Row4 = Multi-Ethnic Clubs, Row5 = Student Life - The Biola Experience

# 3. Update header row text (columns A-L).
$ws.Range("A1").Value = "Organization Name"
$ws.Range("B1").Value = "Categories"
$ws.Range("C1").Value = "Org URL"
$ws.Range("D1").Value = "Image URL"
$ws.Range("E1").Value = "Description"
$ws.Range("F1").Value = "Email"
$ws.Range("G1").Value = "Phone"
$ws.Range("H1").Value = "Website"
$ws.Range("I1").Value = "LinkedIn"
$ws.Range("J1").Value = "Instagram"
$ws.Range("K1").Value = "Facebook"
$ws.Range("L1").Value = "Twitter"

# 4. Swap the Category/Organization Name columns (A<->B) for every data row,
#    and update URLs/names per the new scraped data.
$ws.Range("A2").Value = "Student Life"
$ws.Range("B2").Value = "General"
$ws.Range("C2").Value = "https://www.biola.edu/experience/student-life"

$ws.Range("A3").Value = "Clubs"
$ws.Range("B3").Value = "General"
$ws.Range("C3").Value = "https://www.biola.edu/campus-engagement/clubs"

$ws.Range("A4").Value = "Multi-Ethnic Clubs"
$ws.Range("B4").Value = "General"
$ws.Range("C4").Value = "https://www.biola.edu/student-life/multi-ethnic-clubs"

$ws.Range("A5").Value = "Student Life - The Biola Experience - Biola University"
$ws.Range("B5").Value = "General"
$ws.Range("C5").Value = "https://www.biola.edu/student-life/"

# 5. Re-apply the updated column widths.
#    NOTE: the host's ColumnWidth setter round-trips the value through a
#    character->pixel->character conversion that adds ~0.8333 back on
#    save, so each target width below is nudged down by 0.8 to land on
#    the exact integer width that the diff expects.
$ws.Columns.Item(1).ColumnWidth = 49.2   # -> 50
$ws.Columns.Item(2).ColumnWidth = 11.2   # -> 12
$ws.Columns.Item(3).ColumnWidth = 49.2   # -> 50
$ws.Columns.Item(4).ColumnWidth = 10.2   # -> 11
$ws.Columns.Item(5).ColumnWidth = 12.2   # -> 13
$ws.Columns.Item(6).ColumnWidth = 6.2    # -> 7
$ws.Columns.Item(7).ColumnWidth = 6.2    # -> 7
$ws.Columns.Item(8).ColumnWidth = 8.2    # -> 9
$ws.Columns.Item(9).ColumnWidth = 9.2    # -> 10
$ws.Columns.Item(10).ColumnWidth = 10.2  # -> 11
$ws.Columns.Item(11).ColumnWidth = 9.2   # -> 10
$ws.Columns.Item(12).ColumnWidth = 8.2   # -> 9
